# Sue Sampson.docx - Mod-5 test file update (5/27/2024)
#
# 1. Append a trailing space run to the end of the first ("Sue Sampson -
#    Testme!") paragraph.
# 2. Append a new paragraph: "Test me again…… 5/17/2024" with the date
#    wrapped in gramStart/gramEnd proofing-error markers (as Word's
#    grammar checker would flag it).
# 3. Leave a trailing empty paragraph after that, before the section
#    properties.

$d = $word.ActiveDocument

# --- 1. trailing space on the existing last paragraph -----------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter(" ")

# --- 2 & 3. new paragraph (with proofErr-wrapped date) + blank para ---
# Insert raw OOXML at the (now final) paragraph mark so the gramStart /
# gramEnd proofing marks land exactly around the date run, and the
# paragraph break produces the extra trailing empty paragraph.
$insertPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Test me again&#8230;&#8230; </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>5/17/2024</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($newParaXml)
